$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# New daily records to append (DATA, DESCARTADOS, EM INVESTIGACAO, CONFIRMADOS, EXAMINADOS,
# RECUPERADOS, ATIVOS, HOSPITAL, DOMICILIO, OBITOS)
$data = @(
    @(44364, 12020, 238, 5700, 17958, 5129, 446, 22, 424, 125),
    @(44365, 12135, 214, 5728, 18077, 5134, 466, 23, 443, 128),
    @(44366, 12234, 133, 5772, 18139, 5197, 447, 22, 425, 128),
    @(44367, 12253, 97, 5789, 18139, 5242, 419, 23, 396, 128)
)

$startRow = 351
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    for ($col = 1; $col -le 10; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $values[$col - 1]
    }
}

$lastRow = $startRow + $data.Count - 1

# Update the frozen pane / selection to reflect the new extent of the data.
$ws.Activate()
$ws.Range("A" + ($lastRow - 5)).Select()
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("J" + ($lastRow - 1)).Select()
